$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "1.003"); Excel would
# auto-convert these to numbers unless the cell is forced to Text first. We
# flip the whole D2:D51 range to Text, write the values, then restore the
# Normal style so no residual per-cell style index is left behind (matches
# the unstyled cells in the source workbook).
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"
$ws.Range("D2").Value = "27.105.77"
$ws.Range("D3").Value = "1.894.44"
$ws.Range("D4").Value = "1.003"
$ws.Range("D5").Value = "306.64"
$ws.Range("D7").Value = "0.5203"
$ws.Range("D8").Value = "0.3763"
$ws.Range("D9").Value = "0.07260"
$ws.Range("D11").Value = "0.8999"
$ws.Range("D12").Value = "0.08188"
$ws.Range("D13").Value = "1.970.20"
$ws.Range("D14").Value = "96.35"
$ws.Range("D15").Value = "5.304"
$ws.Range("D17").Value = "0.000008599"
$ws.Range("D19").Value = "1.003"
$ws.Range("D20").Value = "27.142.44"
$ws.Range("D21").Value = "5.078"
$ws.Range("D22").Value = "2.212.64"
$ws.Range("D23").Value = "10.70"
$ws.Range("D24").Value = "6.419"
$ws.Range("D25").Value = "2.318"
$ws.Range("D26").Value = "148.27"
$ws.Range("D27").Value = "18.18"
$ws.Range("D28").Value = "1.735"
$ws.Range("D29").Value = "115.16"
$ws.Range("D30").Value = "4.797"
$ws.Range("D31").Value = "4.858"
$ws.Range("D32").Value = "0.09200"
$ws.Range("D33").Value = "0.05018"
$ws.Range("D34").Value = "0.7928"
$ws.Range("D35").Value = "1.218"
$ws.Range("D36").Value = "3.438"
$ws.Range("D37").Value = "2.965"
$ws.Range("D38").Value = "2.608"
$ws.Range("D39").Value = "0.5715"
$ws.Range("D40").Value = "0.01997"
$ws.Range("D41").Value = "1.074"
$ws.Range("D42").Value = "9.025"
$ws.Range("D43").Value = "6.559"
$ws.Range("D44").Value = "116.38"
$ws.Range("D45").Value = "0.1514"
$ws.Range("D46").Value = "0.4868"
$ws.Range("D47").Value = "1.001"
$ws.Range("D48").Value = "10.07"
$ws.Range("D49").Value = "1.621"
$ws.Range("D50").Value = "38.26"
$ws.Range("D51").Value = "63.63"
$dRange.Style = "Normal"

# Coin name / link / volume columns are plain text already, no coercion risk.
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("E13").Value = "  +3.24%  "
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("E22").Value = "  +2.65%  "
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E31").Value = "  -2.81%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E34").Value = "  -1.81%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E35").Value = "  -2.15%  "
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("E36").Value = "  +2.29%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E38").Value = "  +1.18%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("E44").Value = "  -2.75%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E48").Value = "  -1.86%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("E50").Value = "  +1.53%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("E51").Value = "  -0.35%  "
